$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$Bvals = @(1.270828066841489,1.125245784568278,1.035628213018413,0.9990528910933563,0.9929763059291759,1.035135166862688,1.220680227886874,1.582631768227827,1.847316427812586,1.967443578893779,2.012890627224067,2.003104719275427,1.971183399930737,1.95162506673131,1.839460016149758,1.770577196499971,1.730931476493936,1.717503698548967,1.777912613125807,1.980560637090605,2.112754065083095,2.042223474119169,1.774596408365426,1.484926450598039)
for ($i = 0; $i -lt $Bvals.Length; $i++) {
    $row = $i + 2
    $ws.Range("B$row").Value2 = $Bvals[$i]
}

$Cvals = @(0.3351633210374416,0.2934747437915064,0.267787757299601,0.257298084929829,0.2555549715297047,0.2676463782238159,0.3208081443543165,0.4243217106224506,0.4999020587448513,0.5341791259980937,0.5471433876454626,0.5443520113158229,0.5352460223380717,0.529666274879105,0.4976598071306739,0.4779975674233015,0.4666785400096387,0.462844443856568,0.4800916682034995,0.5379211037891309,0.5756240606892788,0.5555099059781128,0.4791449717453702,0.3963997083005779)
for ($i = 0; $i -lt $Cvals.Length; $i++) {
    $row = $i + 2
    $ws.Range("C$row").Value2 = $Cvals[$i]
}

$Dvals = @(0.3524890677164194,0.3410804870248398,0.334193557090174,0.3314168884784721,0.3309576295423398,0.3341559890702825,0.3485309826342871,0.3776519215137739,0.3996115039553558,0.4097234860183789,0.4135701387348831,0.4127409187954072,0.4100396032604863,0.4083872396279844,0.3989531139084193,0.3931968443194478,0.389897525177247,0.3887824179645634,0.3938084162806206,0.4108325722769735,0.4220605730826605,0.4160587168169343,0.3935318933239103,0.3696745682062499)
for ($i = 0; $i -lt $Dvals.Length; $i++) {
    $row = $i + 2
    $ws.Range("D$row").Value2 = $Dvals[$i]
}

$Fvals = @(1.100842768444082,1.103164224985726,1.105489656491166,1.106663095342185,1.10687156697584,1.105504568305562,1.101456087326554,1.100683343992571,1.10452092266182,1.107231304153089,1.108396978366329,1.108139724156132,1.107324409446591,1.106843165815079,1.104363255028673,1.103089404027656,1.102447480203153,1.102245706039199,1.103215609833242,1.107560101513599,1.111211768851405,1.109188283665972,1.103158270545237,1.100121409913299)
for ($i = 0; $i -lt $Fvals.Length; $i++) {
    $row = $i + 2
    $ws.Range("F$row").Value2 = $Fvals[$i]
}

$Gvals = @(0.002419336053547935,0.002422345133421985,0.002424289857444307,0.00242510684867682,0.002425243991366793,0.002424300776508126,0.002420353468477795,0.002413380134712372,0.002408719812834187,0.002406699222735451,0.00240594829636226,0.002406109389971494,0.002406637158813447,0.002406962283341487,0.002408853857706227,0.002410039690823389,0.002410731112065593,0.002410966825226568,0.002409912488455573,0.002406481754852438,0.00240432246686685,0.002405467357165464,0.002409969966621627,0.002415184956424705)
for ($i = 0; $i -lt $Gvals.Length; $i++) {
    $row = $i + 2
    $ws.Range("G$row").Value2 = $Gvals[$i]
}

$Ivals = @(0.4114338666674886,0.4218043733499783,0.4285951286925673,0.4314685176079429,0.4319520404550428,0.4286334510499241,0.4149215854687558,0.3914047441400239,0.3762012930754715,0.3697392604144447,0.3673578993477822,0.3678678414189314,0.3695420261392286,0.3705760759071488,0.3766327714567339,0.380464927141805,0.3827118008876251,0.3834798794462237,0.380052563656017,0.3690484922582105,0.3622396460761426,0.3658384923498694,0.3802388570390836,0.3974037400132016)
for ($i = 0; $i -lt $Ivals.Length; $i++) {
    $row = $i + 2
    $ws.Range("I$row").Value2 = $Ivals[$i]
}

$Jvals = @(0.3477075488111439,0.3361625255661949,0.3292828434728818,0.3265317684174107,0.3260781193712603,0.3292455292769461,0.3436833767600689,0.3736610294894263,0.3967127096644418,0.40742530651697,0.4115145914469593,0.4106324364329481,0.4077610792678144,0.4060065465576486,0.3960171796216372,0.3899470649045469,0.3864769860748396,0.385305730063294,0.390591034620968,0.4086035792131923,0.4205662243296331,0.4141640736745416,0.3902998345282782,0.3653717823486744)
for ($i = 0; $i -lt $Jvals.Length; $i++) {
    $row = $i + 2
    $ws.Range("J$row").Value2 = $Jvals[$i]
}

$Ovals = @(2.213463880944857,2.236217478565891,2.252292274858519,2.259370367914286,2.260577492119651,2.252385598848477,2.220871683690419,2.17583963973928,2.15307757520219,2.144985989633824,2.142249093695312,2.142823951782361,2.144754255461919,2.145979288129951,2.153652070143096,2.158939958786846,2.162194308690516,2.163332686021363,2.158355007024397,2.144178383735721,2.136821214543545,2.140572692319182,2.158618796393341,2.186216508255541)
for ($i = 0; $i -lt $Ovals.Length; $i++) {
    $row = $i + 2
    $ws.Range("O$row").Value2 = $Ovals[$i]
}
